# Fill in the "e1a" model's error_output column (F) for each prompt-type
# row, and complete the one missing "a19f" (D) value for the "cove" row -
# i.e. finish entering the phase-1 categorization results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - bad prompt
$ws.Range("F4").NumberFormat = "0.00%"
$ws.Range("F4").Value = 0.0025

# Row 5 - good prompt
$ws.Range("F5").NumberFormat = "0%"
$ws.Range("F5").Value = 0

# Row 6 - cot
$ws.Range("F6").NumberFormat = "0%"
$ws.Range("F6").Value = 0

# Row 7 - cove (also finish the previously-blank a19f value)
$ws.Range("D7").Value = 0
$ws.Range("F7").NumberFormat = "0%"
$ws.Range("F7").Value = 0

# Row 8 - one shot
$ws.Range("F8").NumberFormat = "0.00%"
$ws.Range("F8").Value = 0

# Row 9 - few shot
$ws.Range("F9").NumberFormat = "0%"
$ws.Range("F9").Value = 0

# Leave the selection where the author's cursor ended up
$ws.Activate()
$ws.Range("G11").Select() | Out-Null
